$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells hold price text like "27.738.79" or "327.32" that Excel would
# otherwise auto-convert to a Number (losing formatting / exact text). Force
# the cell format to Text before writing so the literal string is preserved,
# exactly matching the original inline-string cell type in the workbook.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.738.79"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.774.62"
$ws.Range("E3").Value = "  +1.40%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.32"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4574"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3584"
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07488"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.82"
$ws.Range("E10").Value = "  -0.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.104"
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.82"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.042"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.211"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.773.95"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.62"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001060"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06447"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.09"
$ws.Range("E21").Value = "  +1.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.804"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.783.06"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.32"
$ws.Range("E24").Value = "  +1.00%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.03"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.28"
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.980.14"
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.184"
$ws.Range("E29").Value = "  +4.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.60"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.102"
$ws.Range("E31").Value = "  +1.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09200"
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.534"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.86"
$ws.Range("E35").Value = "  -0.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02292"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06168"
$ws.Range("E37").Value = "  +2.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2089"
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6322"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.969"
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("E41").Value = "  -1.46%  "
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.823"
$ws.Range("E43").Value = "  +0.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.30"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5919"
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.78"
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06922"
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.138"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.47"
$ws.Range("E51").Value = "  +0.55%  "
